# Daily attendance processing - 2026-01-09 17:08:40
# Swap the order of the two comma-separated names/emails in the
# "Recorded By" column (G) wherever a "System" / "admin@admin.com" entry
# is paired with another recorder, e.g. "System, dnasr281@gmail.com"
# becomes "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Exact "Recorded By" values that need their two parts reversed.
$targets = @(
    "System, dnasr281@gmail.com",
    "admin@admin.com, System",
    "admin@admin.com, dnasr281@gmail.com"
)

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Text

    if ($targets -contains $v) {
        $parts = $v.Split(",")
        $first = $parts[0].Trim()
        $second = $parts[1].Trim()
        $cell.Value = "$second, $first"
    }
}
